{"js": "// Update the division-fact answers in the practice table.\n// Each populated row of the table holds 5 answers (one per column); the\n// table also contains blank spacer rows which are left untouched.\n// We address each target cell explicitly by (row, column) so that the two\n// cells that originally shared identical text (\"65\u00f79=7, 2\") can still be\n// updated to their own distinct replacement values. Using\n// paragraph.insertText(..., Word.InsertLocation.replace) swaps only the\n// text of the existing run, leaving paragraph/run formatting untouched.\n\nconst updates = [\n  { row: 0, col: 0, oldText: \"12\u00f79=1, 3\", newText: \"85\u00f72=42, 1\" },\n  { row: 0, col: 1, oldText: \"10\u00f77=1, 3\", newText: \"87\u00f74=21, 3\" },\n  { row: 0, col: 2, oldText: \"19\u00f76=3, 1\", newText: \"55\u00f74=13, 3\" },\n  { row: 0, col: 3, oldText: \"70\u00f72=35, 0\", newText: \"39\u00f73=13, 0\" },\n  { row: 0, col: 4, oldText: \"18\u00f75=3, 3\", newText: \"49\u00f79=5, 4\" },\n\n  { row: 4, col: 0, oldText: \"40\u00f75=8, 0\", newText: \"64\u00f73=21, 1\" },\n  { row: 4, col: 1, oldText: \"66\u00f79=7, 3\", newText: \"10\u00f78=1, 2\" },\n  { row: 4, col: 2, oldText: \"65\u00f77=9, 2\", newText: \"53\u00f78=6, 5\" },\n  { row: 4, col: 3, oldText: \"56\u00f73=18, 2\", newText: \"99\u00f79=11, 0\" },\n  { row: 4, col: 4, oldText: \"95\u00f76=15, 5\", newText: \"28\u00f78=3, 4\" },\n\n  { row: 8, col: 0, oldText: \"65\u00f79=7, 2\", newText: \"33\u00f73=11, 0\" },\n  { row: 8, col: 1, oldText: \"70\u00f77=10, 0\", newText: \"73\u00f73=24, 1\" },\n  { row: 8, col: 2, oldText: \"85\u00f77=12, 1\", newText: \"30\u00f76=5, 0\" },\n  { row: 8, col: 3, oldText: \"22\u00f79=2, 4\", newText: \"23\u00f79=2, 5\" },\n  { row: 8, col: 4, oldText: \"65\u00f79=7, 2\", newText: \"23\u00f79=2, 5\" },\n\n  { row: 12, col: 0, oldText: \"50\u00f77=7, 1\", newText: \"80\u00f73=26, 2\" },\n  { row: 12, col: 1, oldText: \"71\u00f79=7, 8\", newText: \"40\u00f75=8, 0\" },\n  { row: 12, col: 2, oldText: \"86\u00f79=9, 5\", newText: \"13\u00f75=2, 3\" },\n  { row: 12, col: 3, oldText: \"82\u00f74=20, 2\", newText: \"61\u00f74=15, 1\" },\n  { row: 12, col: 4, oldText: \"43\u00f74=10, 3\", newText: \"66\u00f74=16, 2\" },\n\n  { row: 16, col: 0, oldText: \"21\u00f77=3, 0\", newText: \"47\u00f72=23, 1\" },\n  { row: 16, col: 1, oldText: \"24\u00f78=3, 0\", newText: \"82\u00f72=41, 0\" },\n  { row: 16, col: 2, oldText: \"80\u00f76=13, 2\", newText: \"86\u00f76=14, 2\" },\n  { row: 16, col: 3, oldText: \"69\u00f79=7, 6\", newText: \"56\u00f74=14, 0\" },\n  { row: 16, col: 4, oldText: \"90\u00f76=15, 0\", newText: \"58\u00f75=11, 3\" },\n];\n\nconst table = context.document.body.tables.getFirstOrNullObject();\ntable.load(\"rowCount\");\nawait context.sync();\n\nif (table.isNullObject) {\n  throw new Error(\"Expected a table in the document body, but none was found.\");\n}\n\n// Grab the first paragraph of each target cell up front.\nconst paragraphs = updates.map(({ row, col }) => {\n  const cell = table.getCell(row, col);\n  const para = cell.body.paragraphs.getFirst();\n  para.load(\"text\");\n  return para;\n});\nawait context.sync();\n\n// Verify we are editing the expected cells, then replace their text in\n// place (this keeps the existing run/paragraph formatting intact).\nupdates.forEach(({ oldText, newText }, i) => {\n  const para = paragraphs[i];\n  if (para.text !== oldText) {\n    throw new Error(\n      `Unexpected cell text at index ${i}: expected \"${oldText}\" but found \"${para.text}\"`\n    );\n  }\n  para.insertText(newText, Word.InsertLocation.replace);\n});\n\nawait context.sync();\n", "ps1": "# Update the division-fact answers in the practice table.\n# Each populated row of the table holds 5 answers (one per column); the\n# table also contains blank spacer rows which are left untouched.\n# We address each target cell explicitly by (row, column) so that the two\n# cells that originally shared identical text (\"65\u00f79=7, 2\") can still be\n# updated to their own distinct replacement values.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$updates = @(\n    @{ Row = 1;  Col = 1; Old = \"12\u00f79=1, 3\";  New = \"85\u00f72=42, 1\" },\n    @{ Row = 1;  Col = 2; Old = \"10\u00f77=1, 3\";  New = \"87\u00f74=21, 3\" },\n    @{ Row = 1;  Col = 3; Old = \"19\u00f76=3, 1\";  New = \"55\u00f74=13, 3\" },\n    @{ Row = 1;  Col = 4; Old = \"70\u00f72=35, 0\"; New = \"39\u00f73=13, 0\" },\n    @{ Row = 1;  Col = 5; Old = \"18\u00f75=3, 3\";  New = \"49\u00f79=5, 4\" },\n\n    @{ Row = 5;  Col = 1; Old = \"40\u00f75=8, 0\";   New = \"64\u00f73=21, 1\" },\n    @{ Row = 5;  Col = 2; Old = \"66\u00f79=7, 3\";   New = \"10\u00f78=1, 2\" },\n    @{ Row = 5;  Col = 3; Old = \"65\u00f77=9, 2\";   New = \"53\u00f78=6, 5\" },\n    @{ Row = 5;  Col = 4; Old = \"56\u00f73=18, 2\";  New = \"99\u00f79=11, 0\" },\n    @{ Row = 5;  Col = 5; Old = \"95\u00f76=15, 5\";  New = \"28\u00f78=3, 4\" },\n\n    @{ Row = 9;  Col = 1; Old = \"65\u00f79=7, 2\";  New = \"33\u00f73=11, 0\" },\n    @{ Row = 9;  Col = 2; Old = \"70\u00f77=10, 0\"; New = \"73\u00f73=24, 1\" },\n    @{ Row = 9;  Col = 3; Old = \"85\u00f77=12, 1\"; New = \"30\u00f76=5, 0\" },\n    @{ Row = 9;  Col = 4; Old = \"22\u00f79=2, 4\";  New = \"23\u00f79=2, 5\" },\n    @{ Row = 9;  Col = 5; Old = \"65\u00f79=7, 2\";  New = \"23\u00f79=2, 5\" },\n\n    @{ Row = 13; Col = 1; Old = \"50\u00f77=7, 1\";  New = \"80\u00f73=26, 2\" },\n    @{ Row = 13; Col = 2; Old = \"71\u00f79=7, 8\";  New = \"40\u00f75=8, 0\" },\n    @{ Row = 13; Col = 3; Old = \"86\u00f79=9, 5\";  New = \"13\u00f75=2, 3\" },\n    @{ Row = 13; Col = 4; Old = \"82\u00f74=20, 2\"; New = \"61\u00f74=15, 1\" },\n    @{ Row = 13; Col = 5; Old = \"43\u00f74=10, 3\"; New = \"66\u00f74=16, 2\" },\n\n    @{ Row = 17; Col = 1; Old = \"21\u00f77=3, 0\"; New = \"47\u00f72=23, 1\" },\n    @{ Row = 17; Col = 2; Old = \"24\u00f78=3, 0\"; New = \"82\u00f72=41, 0\" },\n    @{ Row = 17; Col = 3; Old = \"80\u00f76=13, 2\"; New = \"86\u00f76=14, 2\" },\n    @{ Row = 17; Col = 4; Old = \"69\u00f79=7, 6\"; New = \"56\u00f74=14, 0\" },\n    @{ Row = 17; Col = 5; Old = \"90\u00f76=15, 0\"; New = \"58\u00f75=11, 3\" }\n)\n\nforeach ($u in $updates) {\n    $cell = $t.Cell($u.Row, $u.Col)\n    # Cell.Range.Text includes the trailing cell-mark characters, so trim\n    # those off before comparing with the expected current value.\n    $current = $cell.Range.Text.TrimEnd([char]13, [char]7)\n    if ($current -ne $u.Old) {\n        throw \"Unexpected text in cell ($($u.Row),$($u.Col)): expected '$($u.Old)' but found '$current'\"\n    }\n    $cell.Range.Text = $u.New\n}\n"}
